$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-72 down to 70-73
$ws.Rows.Item(69).Insert()

# Copy formatting of the new row 69 from row 70 (which now holds what used to be row 69)
$ws.Range("A70:R70").Copy()
$ws.Range("A69:R69").PasteSpecial()

# Fill in the new row 69 with its data (constant columns copied from the surrounding rows)
$ws.Cells.Item(69, 1).Value = 11
$ws.Cells.Item(69, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(69, 3).Value = "Bíobío"
$ws.Cells.Item(69, 4).Value = 44516
$ws.Cells.Item(69, 5).Value = 8
$ws.Cells.Item(69, 6).Value = 100112032
$ws.Cells.Item(69, 7).Value = "Zapallo italiano"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 350
$ws.Cells.Item(69, 11).Value = 4500
$ws.Cells.Item(69, 12).Value = 5000
$ws.Cells.Item(69, 13).Value = 4786
$ws.Cells.Item(69, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(69, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(69, 16).Value = 80
$ws.Cells.Item(69, 17).Value = 60
$ws.Cells.Item(69, 18).Value = "Hortaliza"
